$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference, new value, whether to force Text number format
# (columns D/E/G hold numeric-looking strings - e.g. "328.06", "-0.68%", "15" -
#  these must be forced to Text format so Excel keeps them as literal strings
#  instead of auto-converting to numbers/percentages).
$updates = @(
    @{ Cell = "D2"; Value = "328.06"; ForceText = $true }
    @{ Cell = "E2"; Value = "-0.68%"; ForceText = $true }
    @{ Cell = "G2"; Value = "15"; ForceText = $true }
    @{ Cell = "D3"; Value = "43.94"; ForceText = $true }
    @{ Cell = "E3"; Value = "1.45%"; ForceText = $true }
    @{ Cell = "G3"; Value = "15"; ForceText = $true }
    @{ Cell = "D4"; Value = "5.566"; ForceText = $true }
    @{ Cell = "E4"; Value = "-0.71%"; ForceText = $true }
    @{ Cell = "G4"; Value = "15"; ForceText = $true }
    @{ Cell = "D5"; Value = "0.08047"; ForceText = $true }
    @{ Cell = "E5"; Value = "-1.95%"; ForceText = $true }
    @{ Cell = "G5"; Value = "15"; ForceText = $true }
    @{ Cell = "E6"; Value = "-0.62%"; ForceText = $true }
    @{ Cell = "G6"; Value = "15"; ForceText = $true }
    @{ Cell = "D7"; Value = "4.270"; ForceText = $true }
    @{ Cell = "E7"; Value = "-2.74%"; ForceText = $true }
    @{ Cell = "G7"; Value = "15"; ForceText = $true }
    @{ Cell = "D8"; Value = "0.9448"; ForceText = $true }
    @{ Cell = "E8"; Value = "0.01%"; ForceText = $true }
    @{ Cell = "G8"; Value = "15"; ForceText = $true }
    @{ Cell = "D9"; Value = "2.520"; ForceText = $true }
    @{ Cell = "E9"; Value = "-10.21%"; ForceText = $true }
    @{ Cell = "G9"; Value = "15"; ForceText = $true }
    @{ Cell = "D10"; Value = "0.1173"; ForceText = $true }
    @{ Cell = "E10"; Value = "-2.03%"; ForceText = $true }
    @{ Cell = "G10"; Value = "15"; ForceText = $true }
    @{ Cell = "D11"; Value = "0.1840"; ForceText = $true }
    @{ Cell = "E11"; Value = "-4.24%"; ForceText = $true }
    @{ Cell = "G11"; Value = "15"; ForceText = $true }
    @{ Cell = "D12"; Value = "0.09659"; ForceText = $true }
    @{ Cell = "E12"; Value = "-1.24%"; ForceText = $true }
    @{ Cell = "G12"; Value = "15"; ForceText = $true }
    @{ Cell = "D13"; Value = "0.04365"; ForceText = $true }
    @{ Cell = "E13"; Value = "0.35%"; ForceText = $true }
    @{ Cell = "G13"; Value = "15"; ForceText = $true }
    @{ Cell = "D14"; Value = "0.1066"; ForceText = $true }
    @{ Cell = "E14"; Value = "-0.39%"; ForceText = $true }
    @{ Cell = "G14"; Value = "15"; ForceText = $true }
    @{ Cell = "D15"; Value = "0.001276"; ForceText = $true }
    @{ Cell = "E15"; Value = "-0.88%"; ForceText = $true }
    @{ Cell = "G15"; Value = "15"; ForceText = $true }
    @{ Cell = "D16"; Value = "0.005982"; ForceText = $true }
    @{ Cell = "E16"; Value = "1.03%"; ForceText = $true }
    @{ Cell = "G16"; Value = "15"; ForceText = $true }
    @{ Cell = "D17"; Value = "3.397"; ForceText = $true }
    @{ Cell = "E17"; Value = "-3.52%"; ForceText = $true }
    @{ Cell = "G17"; Value = "15"; ForceText = $true }
    @{ Cell = "D18"; Value = "0.3450"; ForceText = $true }
    @{ Cell = "E18"; Value = "-2.44%"; ForceText = $true }
    @{ Cell = "G18"; Value = "15"; ForceText = $true }
    @{ Cell = "D19"; Value = "10.30"; ForceText = $true }
    @{ Cell = "E19"; Value = "18.27%"; ForceText = $true }
    @{ Cell = "G19"; Value = "15"; ForceText = $true }
    @{ Cell = "D20"; Value = "0.1380"; ForceText = $true }
    @{ Cell = "E20"; Value = "0.79%"; ForceText = $true }
    @{ Cell = "G20"; Value = "15"; ForceText = $true }
    @{ Cell = "D21"; Value = "0.2504"; ForceText = $true }
    @{ Cell = "E21"; Value = "-0.76%"; ForceText = $true }
    @{ Cell = "G21"; Value = "15"; ForceText = $true }
    @{ Cell = "D22"; Value = "0.04190"; ForceText = $true }
    @{ Cell = "E22"; Value = "-4.83%"; ForceText = $true }
    @{ Cell = "G22"; Value = "15"; ForceText = $true }
    @{ Cell = "D23"; Value = "0.001246"; ForceText = $true }
    @{ Cell = "E23"; Value = "0.25%"; ForceText = $true }
    @{ Cell = "G23"; Value = "15"; ForceText = $true }
    @{ Cell = "D24"; Value = "0.004280"; ForceText = $true }
    @{ Cell = "E24"; Value = "-0.80%"; ForceText = $true }
    @{ Cell = "G24"; Value = "15"; ForceText = $true }
    @{ Cell = "D25"; Value = "0.0001260"; ForceText = $true }
    @{ Cell = "E25"; Value = "1.88%"; ForceText = $true }
    @{ Cell = "G25"; Value = "15"; ForceText = $true }
    @{ Cell = "D26"; Value = "0.0003989"; ForceText = $true }
    @{ Cell = "E26"; Value = "-0.54%"; ForceText = $true }
    @{ Cell = "G26"; Value = "15"; ForceText = $true }
    @{ Cell = "G27"; Value = "15"; ForceText = $true }
    @{ Cell = "G28"; Value = "15"; ForceText = $true }
    @{ Cell = "G29"; Value = "15"; ForceText = $true }
    @{ Cell = "G30"; Value = "15"; ForceText = $true }
    @{ Cell = "G31"; Value = "15"; ForceText = $true }
    @{ Cell = "G32"; Value = "15"; ForceText = $true }
    @{ Cell = "G33"; Value = "15"; ForceText = $true }
    @{ Cell = "G34"; Value = "15"; ForceText = $true }
    @{ Cell = "G35"; Value = "15"; ForceText = $true }
    @{ Cell = "G36"; Value = "15"; ForceText = $true }
    @{ Cell = "G37"; Value = "15"; ForceText = $true }
    @{ Cell = "D38"; Value = "0.02650"; ForceText = $true }
    @{ Cell = "E38"; Value = "-4.64%"; ForceText = $true }
    @{ Cell = "G38"; Value = "15"; ForceText = $true }
    @{ Cell = "D39"; Value = "0.05506"; ForceText = $true }
    @{ Cell = "E39"; Value = "-4.15%"; ForceText = $true }
    @{ Cell = "G39"; Value = "15"; ForceText = $true }
    @{ Cell = "D40"; Value = "0.007561"; ForceText = $true }
    @{ Cell = "E40"; Value = "-4.78%"; ForceText = $true }
    @{ Cell = "G40"; Value = "15"; ForceText = $true }
    @{ Cell = "D41"; Value = "0.1393"; ForceText = $true }
    @{ Cell = "E41"; Value = "-1.84%"; ForceText = $true }
    @{ Cell = "G41"; Value = "15"; ForceText = $true }
    @{ Cell = "D42"; Value = "0.007903"; ForceText = $true }
    @{ Cell = "E42"; Value = "-18.91%"; ForceText = $true }
    @{ Cell = "G42"; Value = "15"; ForceText = $true }
    @{ Cell = "D43"; Value = "0.002000"; ForceText = $true }
    @{ Cell = "E43"; Value = "-3.38%"; ForceText = $true }
    @{ Cell = "G43"; Value = "15"; ForceText = $true }
    @{ Cell = "D44"; Value = "0.008825"; ForceText = $true }
    @{ Cell = "E44"; Value = "-8.85%"; ForceText = $true }
    @{ Cell = "G44"; Value = "15"; ForceText = $true }
    @{ Cell = "D45"; Value = "0.00006866"; ForceText = $true }
    @{ Cell = "E45"; Value = "-9.76%"; ForceText = $true }
    @{ Cell = "G45"; Value = "15"; ForceText = $true }
    @{ Cell = "D46"; Value = "0.00000000750"; ForceText = $true }
    @{ Cell = "E46"; Value = "-0.54%"; ForceText = $true }
    @{ Cell = "G46"; Value = "15"; ForceText = $true }
    @{ Cell = "B47"; Value = "BOLO"; ForceText = $false }
    @{ Cell = "C47"; Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"; ForceText = $false }
    @{ Cell = "D47"; Value = "0.003372"; ForceText = $true }
    @{ Cell = "E47"; Value = "-2.27%"; ForceText = $true }
    @{ Cell = "G47"; Value = "15"; ForceText = $true }
    @{ Cell = "B48"; Value = "CoinbaseStockToken"; ForceText = $false }
    @{ Cell = "C48"; Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"; ForceText = $false }
    @{ Cell = "D48"; Value = "0.002270"; ForceText = $true }
    @{ Cell = "E48"; Value = "-0.54%"; ForceText = $true }
    @{ Cell = "G48"; Value = "15"; ForceText = $true }
    @{ Cell = "D49"; Value = "0.00002100"; ForceText = $true }
    @{ Cell = "E49"; Value = "-0.54%"; ForceText = $true }
    @{ Cell = "G49"; Value = "15"; ForceText = $true }
    @{ Cell = "D50"; Value = "0.0002000"; ForceText = $true }
    @{ Cell = "E50"; Value = "-0.54%"; ForceText = $true }
    @{ Cell = "G50"; Value = "15"; ForceText = $true }
    @{ Cell = "G51"; Value = "15"; ForceText = $true }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $u.Value
}
